$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "30.168.18"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.838.71"
$ws.Range("E3").Value = "  -1.46%  "
Set-TextValue $ws.Range("D4") "1.0000"
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "233.08"
$ws.Range("E5").Value = "  -0.91%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.07%  "
Set-TextValue $ws.Range("D7") "0.4673"
$ws.Range("E7").Value = "  -3.08%  "
Set-TextValue $ws.Range("D8") "0.2701"
$ws.Range("E8").Value = "  -3.25%  "
Set-TextValue $ws.Range("D9") "0.06269"
$ws.Range("E9").Value = "  -3.58%  "
$ws.Range("D10").Value = "1.838.33"
$ws.Range("E10").Value = "  -2.58%  "
Set-TextValue $ws.Range("D11") "0.07413"
$ws.Range("E11").Value = "  -0.31%  "
Set-TextValue $ws.Range("D12") "16.04"
$ws.Range("E12").Value = "  -1.55%  "
Set-TextValue $ws.Range("D13") "4.924"
$ws.Range("E13").Value = "  -2.87%  "
Set-TextValue $ws.Range("D14") "83.71"
$ws.Range("E14").Value = "  -3.98%  "
Set-TextValue $ws.Range("D15") "0.6197"
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").Value = "30.082.39"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("E17").Value = "  -0.03%  "
Set-TextValue $ws.Range("D18") "228.08"
$ws.Range("E18").Value = "  -2.20%  "
Set-TextValue $ws.Range("D19") "0.000007274"
$ws.Range("E19").Value = "  -3.19%  "
Set-TextValue $ws.Range("D20") "12.33"
$ws.Range("E20").Value = "  -5.05%  "
Set-TextValue $ws.Range("D21") "0.9985"
$ws.Range("E21").Value = "  -0.26%  "
Set-TextValue $ws.Range("D22") "4.865"
$ws.Range("E22").Value = "  -5.55%  "
Set-TextValue $ws.Range("D23") "5.835"
$ws.Range("E23").Value = "  -4.21%  "
Set-TextValue $ws.Range("D24") "9.190"
$ws.Range("E24").Value = "  -1.52%  "
Set-TextValue $ws.Range("D25") "164.37"
$ws.Range("E25").Value = "  -2.06%  "
Set-TextValue $ws.Range("D26") "17.82"
$ws.Range("E26").Value = "  -2.95%  "
Set-TextValue $ws.Range("D27") "1.884"
$ws.Range("E27").Value = "  -1.99%  "
Set-TextValue $ws.Range("D28") "0.1027"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.46%  "
Set-TextValue $ws.Range("D30") "4.065"
$ws.Range("E30").Value = "  -4.70%  "
Set-TextValue $ws.Range("D31") "3.788"
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("E32").Value = "  -3.69%  "
Set-TextValue $ws.Range("D33") "1.136"
$ws.Range("E33").Value = "  -3.28%  "
Set-TextValue $ws.Range("D34") "0.7074"
$ws.Range("E34").Value = "  -4.64%  "
Set-TextValue $ws.Range("D35") "2.691"
$ws.Range("E35").Value = "  -0.74%  "
Set-TextValue $ws.Range("D36") "0.01888"
$ws.Range("E36").Value = "  -2.01%  "
Set-TextValue $ws.Range("D37") "2.646"
$ws.Range("E37").Value = "  +0.32%  "
Set-TextValue $ws.Range("D38") "0.8931"
$ws.Range("E38").Value = "  -3.20%  "
Set-TextValue $ws.Range("D39") "1.922"
$ws.Range("E39").Value = "  -6.41%  "
Set-TextValue $ws.Range("D40") "104.28"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("E41").Value = "  +0.52%  "
Set-TextValue $ws.Range("D42") "5.528"
$ws.Range("E42").Value = "  -0.34%  "
Set-TextValue $ws.Range("D43") "0.3994"
$ws.Range("E43").Value = "  -4.75%  "
Set-TextValue $ws.Range("D44") "6.988"
$ws.Range("E44").Value = "  -3.54%  "
Set-TextValue $ws.Range("D45") "0.1190"
$ws.Range("E45").Value = "  -3.31%  "
Set-TextValue $ws.Range("D46") "59.59"
$ws.Range("E46").Value = "  -3.69%  "
Set-TextValue $ws.Range("D47") "8.507"
$ws.Range("E47").Value = "  -4.76%  "
Set-TextValue $ws.Range("D48") "32.63"
$ws.Range("E48").Value = "  -2.90%  "
Set-TextValue $ws.Range("D49") "0.05505"
$ws.Range("E49").Value = "  -2.55%  "
Set-TextValue $ws.Range("D50") "1.354"
$ws.Range("E50").Value = "  -5.96%  "
Set-TextValue $ws.Range("D51") "0.3615"
$ws.Range("E51").Value = "  -5.06%  "
